# TestData.xlsx: swap in generic "foo author" / "foo description" sample
# text (was the real author's name + per-sheet blurb), and flip which
# sheet/selection is active.

$wb = $excel.ActiveWorkbook

$wsFoo = $wb.Worksheets.Item("Foo")
$wsBar = $wb.Worksheets.Item("Bar")

# --- Foo sheet: update author/description text ---
# (set description before author so the shared-string table gets the
# same insertion order as the target file)
$wsFoo.Range("B2").Value = "foo description"
$wsFoo.Range("B1").Value = "foo author"

# --- Bar sheet: same text updates ---
$wsBar.Range("B2").Value = "foo description"
$wsBar.Range("B1").Value = "foo author"

# --- Selections ---
# Foo: was B10 selected/active; now B1:B2 selected with B2 active.
$wsFoo.Range("B1").Resize(2, 1).Select()

# Bar: was C19 selected/active; now D3 selected/active.
$wsBar.Range("D3").Select()

# --- Active sheet/tab flips from Foo to Bar ---
$wsBar.Activate()
